$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.145737409591675
$ws.Range("B1").Value = 2.207284212112427
$ws.Range("C1").Value = 2.986399412155151
$ws.Range("D1").Value = 1.480307579040527
$ws.Range("E1").Value = 1.037962317466736
